# Export data to Excel and handle pagination:
# Append a new selector row ("suivant" / next-page button xpath) to the
# selectors table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "suivant"
$ws.Range("B11").Value = "/html/body/div[3]/div/section/div[1]/div[1]/div[2]/div[2]/div[2]/nav/div[3]/button"

$ws.Range("A11:B11").Style = $ws.Range("A2:B2").Style
